# Applies the changes described by the commit:
#  - Update the "Date" metadata value on the "Metadata" sheet
#  - Re-order the two "Mapping" columns (and their header/body values) on the
#    "Elements" sheet so the "Spécification métier" mapping column comes
#    before the "RIM Mapping" column
#  - Swap the corresponding column widths to match the new column order

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- 1. Update the Date value on the Metadata sheet -----------------------
$wsMetadata.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Swap the "Mapping: RIM Mapping" (AK) and
#        "Mapping: Spécification métier vers l'extension ROR
#        TelecomConfidentialityLevel" (AL) columns on the Elements sheet ---

$rimMapping = "Mapping: RIM Mapping"
$businessMapping = "Mapping: Spécification métier vers l'extension ROR TelecomConfidentialityLevel"

# Capture the current (pre-swap) values for each row of the two columns
$akValues = @{}
$alValues = @{}
for ($row = 1; $row -le 6; $row++) {
    $akValues[$row] = $wsElements.Range("AK" + $row).Value2
    $alValues[$row] = $wsElements.Range("AL" + $row).Value2
}

# Header row: swap the mapping titles explicitly so the text matches exactly
$wsElements.Range("AK1").Value = $businessMapping
$wsElements.Range("AL1").Value = $rimMapping

# Data rows: swap whatever value each row held between the two columns
for ($row = 2; $row -le 6; $row++) {
    $wsElements.Range("AK" + $row).Value = $alValues[$row]
    $wsElements.Range("AL" + $row).Value = $akValues[$row]
}

# --- 3. Swap the column widths of AK (37) and AL (38) ----------------------
# (Target raw OOXML widths are 84.40625 and 24.98046875; ColumnWidth is
# expressed in characters and gets quantized internally, so these inputs are
# chosen to land in the closest achievable width bucket.)
$wsElements.Columns.Item(37).ColumnWidth = 83.5
$wsElements.Columns.Item(38).ColumnWidth = 24.16666667
